$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list with latest prices / 1h volume changes.
# Some coins moved rank (rows 31/32, 35-37, 47/48 swapped identity),
# so B (Coin), C (Link), D (Price) and E (Volume(1h)) are rewritten
# per-row to match the refreshed data pull.
# NumberFormat is forced to text ("@") before each write so that
# numeric-looking price strings (e.g. "1.00") are stored as text,
# matching the original inline-string cell content/type.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.248.36'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.25%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.184.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.25%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.38'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.41'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.40%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.12'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.99%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.178.25'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.09%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.749'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.50%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.205'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.40%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.16%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.79%  '

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.29%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.203.28'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.19%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.765.01'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.93%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.191.13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.75'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.43%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.25'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +9.66%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.93'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.74%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000206'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '444.83'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.83'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.69%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.39'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.29%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.16'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.83%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.339.52'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.29%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.133'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +51.25%  '

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.233'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +20.45%  '

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +7.93%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.45'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.169'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +11.35%  '

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.986'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.56%  '

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.81'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +10.04%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.57'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.52%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '514.08'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.19%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +7.79%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.95'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.14%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.457'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.81'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +10.83%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.44'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.95%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.15'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.03%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.729'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +7.27%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.23'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.53%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.94'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.20%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.39'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +5.91%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.46'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.89%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.17'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.67%  '
